$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vp_sku_list")

# New SKU row appended at the bottom of the list (A87), matching the
# formatting of A81 (the other "pasted" style row in this column),
# including that row's taller row height.
$ws.Range("A81").Copy()
$ws.Range("A87").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A87").Value = 10025976
$ws.Rows.Item(87).RowHeight = $ws.Rows.Item(81).RowHeight

$ws.Range("A87").Select()
